$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

$changed = 0
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $orig = $cell.Value2
    if ([string]::IsNullOrEmpty($orig)) { continue }

    $parts = $orig -split ',\s*'
    $others = @()
    $systems = @()
    foreach ($p in $parts) {
        if ($p.ToLower() -eq "system") {
            $systems += $p
        } else {
            $others += $p
        }
    }

    if ($systems.Count -gt 0 -and $others.Count -gt 0) {
        $newVal = ($others + $systems) -join ", "
        if ($newVal -ne $orig) {
            $cell.Value = $newVal
            $changed++
        }
    }
}

Write-Host "Changed rows: $changed"
